# Apply the "Add files via upload" edit to ParameterReferenceListList.xlsx
#
# Summary of changes (per the OOXML diff):
#   - sheet1!G39 changes from 30 to 20
#   - a new row 44 is appended with:
#       C44 = "Para.durThr;"
#       D44 = "N/A"
#       E44 = "N/A"
#       F44 = 0.99
#       G44 = 0.99
#       I44 = "The inclusion of this is essential for the morphology detector."
#     (H44 is left empty, matching the diff)
#   - two new shared strings are introduced, in order: "Para.durThr;" then
#     "The inclusion of this is essential for the morphology detector."
#   - the sheet's dimension grows from A1:I43 to A1:I44 (handled automatically
#     by the engine once the new row is populated)
#   - the view's selection moves from H2 to G40, and the view scrolls so that
#     A19 is the top-left visible cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edits -----------------------------------------------------------

# G39: 30 -> 20
$ws.Range("G39").Value = 20

# New row 44
$ws.Range("C44").Value = "Para.durThr;"
$ws.Range("D44").Value = "N/A"
$ws.Range("E44").Value = "N/A"
$ws.Range("F44").Value = 0.99
$ws.Range("G44").Value = 0.99
$ws.Range("I44").Value = "The inclusion of this is essential for the morphology detector."

# --- view state -------------------------------------------------------------
# Select G40 (new active cell) and scroll the window so A19 is the top-left
# visible cell, matching the author's saved view.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("G40").Select()
